$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "...posición de una pieza. Para el " -> "...Para la " (only the one
#    occurrence right before the existing `evento "drop"` run - anchored
#    with enough surrounding context to make the match unique).
# ---------------------------------------------------------------------
$rngFind = $d.Content
$null = $rngFind.Find.Execute(
    "posición de una pieza. Para el ", $true, $false, $false, $false, $false,
    $true, 1, $false, "posición de una pieza. Para la ", 2)

# ---------------------------------------------------------------------
# 2. Locate the insertion point: right after "...Para la " and right
#    before the existing run that reads evento "drop".
# ---------------------------------------------------------------------
$rngFind2 = $d.Content
$null = $rngFind2.Find.Execute(
    "posición de una pieza. Para la ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

$insertStart = $rngFind2.End

# ---------------------------------------------------------------------
# 3. Insert the whole new block of text in one shot (as plain text, so it
#    merges cleanly with the preceding plain run), then go back and apply
#    italic / green-color formatting to the specific sub-ranges that need
#    it. This avoids ever having to "turn off" formatting (which would
#    leave stray explicit w:val="0" / w:val="auto" attributes).
# ---------------------------------------------------------------------
$newText = "modificación del estilo de una pieza cuando se arrastra mediante el evento “dragstart” he decidido añadir transparencia a la imagen que se restablecerá en el evento “drop”, y para cada pieza donde nos situamos con el evento “dragover” he añadido un resaltado azul al elemento que se desactivará mediante el evento “dragleave”. Para el "

$ip = $d.Range($insertStart, $insertStart)
$ip.InsertAfter($newText)

# Green color used for the event-name runs (RGB 38761d -> BGR long).
$greenColor = 1930808

function Style-Segment($relStart, $relEnd, $italic, $color) {
    $r = $d.Range($insertStart + $relStart, $insertStart + $relEnd)
    if ($italic) {
        $r.Font.Italic = 1
    }
    if ($color -ne $null) {
        $r.Font.Color = $color
    }
}

# "modificación del estilo" - italic only
Style-Segment 0 23 $true $null

# evento "dragstart" - italic + green
Style-Segment 68 86 $true $greenColor

# evento "drop" - italic + green
Style-Segment 158 171 $true $greenColor

# evento "dragover" - italic + green
Style-Segment 217 234 $true $greenColor

# evento "dragleave" - italic + green
Style-Segment 307 325 $true $greenColor

Write-Output "Done"
